# Auto-generated Excel COM-interop script
# Applies scheduled-runner profit recalculations to Tiamat_Profits sheets
$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H129").Value = 803.59375
$ws.Range("J129").Value = 1007.7368
$ws.Range("L129").Value = 3023.2104
$ws.Range("N129").Value = -13023.2104
$ws.Range("H132").Value = 164121.42
$ws.Range("I132").Value = 2805.7856
$ws.Range("J132").Value = 1669734
$ws.Range("K132").Value = 8417.356800000001
$ws.Range("L132").Value = 5009202
$ws.Range("M132").Value = -5887.356800000001
$ws.Range("N132").Value = -5014262
$ws.Range("H137").Value = 4408.5625
$ws.Range("I137").Value = 859.3
$ws.Range("J137").Value = 6021.864
$ws.Range("K137").Value = 2577.9
$ws.Range("L137").Value = 18065.592
$ws.Range("M137").Value = -27.89999999999964
$ws.Range("N137").Value = -23165.592
$ws.Range("H138").Value = 2651.9307
$ws.Range("I138").Value = 1861.7307
$ws.Range("J138").Value = 3098.5652
$ws.Range("K138").Value = 5585.1921
$ws.Range("L138").Value = 9295.695599999999
$ws.Range("M138").Value = -445.1921000000002
$ws.Range("N138").Value = -19575.6956

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H24").Value = 29800
$ws.Range("J24").Value = 29800
$ws.Range("L24").Value = 29800
$ws.Range("N24").Value = -30548
$ws.Range("H61").Value = 3784.4614
$ws.Range("I61").Value = 3588.6667
$ws.Range("J61").Value = 4225
$ws.Range("K61").Value = 3588.6667
$ws.Range("L61").Value = 4225
$ws.Range("M61").Value = -3376.6667
$ws.Range("N61").Value = -4649
$ws.Range("H74").Value = 30827.883
$ws.Range("I74").Value = 42947.293
$ws.Range("J74").Value = 1741.3
$ws.Range("K74").Value = 42947.293
$ws.Range("L74").Value = 1741.3
$ws.Range("M74").Value = -42073.293
$ws.Range("N74").Value = -3489.3
$ws.Range("H77").Value = 30827.883
$ws.Range("I77").Value = 42947.293
$ws.Range("J77").Value = 1741.3
$ws.Range("K77").Value = 214736.465
$ws.Range("L77").Value = 8706.5
$ws.Range("M77").Value = -210368.465
$ws.Range("N77").Value = -17442.5
$ws.Range("H100").Value = 29800
$ws.Range("J100").Value = 29800
$ws.Range("L100").Value = 29800
$ws.Range("N100").Value = -31964
$ws.Range("H132").Value = 198386.3
$ws.Range("I132").Value = 31185.35
$ws.Range("J132").Value = 591800.3
$ws.Range("K132").Value = 93556.04999999999
$ws.Range("L132").Value = 1775400.9
$ws.Range("M132").Value = -91026.04999999999
$ws.Range("N132").Value = -1780460.9
$ws.Range("H136").Value = 3784.4614
$ws.Range("I136").Value = 3588.6667
$ws.Range("J136").Value = 4225
$ws.Range("K136").Value = 10766.0001
$ws.Range("L136").Value = 12675
$ws.Range("M136").Value = -8216.000100000001
$ws.Range("N136").Value = -17775

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 335775.28
$ws.Range("I86").Value = 1861.5
$ws.Range("J86").Value = 541260.7
$ws.Range("K86").Value = 1861.5
$ws.Range("L86").Value = 541260.7
$ws.Range("M86").Value = -738.5
$ws.Range("N86").Value = -543506.7
$ws.Range("H89").Value = 335775.28
$ws.Range("I89").Value = 1861.5
$ws.Range("J89").Value = 541260.7
$ws.Range("K89").Value = 9307.5
$ws.Range("L89").Value = 2706303.5
$ws.Range("M89").Value = -3691.5
$ws.Range("N89").Value = -2717535.5
$ws.Range("H134").Value = 22751868
$ws.Range("I134").Value = 1806.3158
$ws.Range("J134").Value = 166835580
$ws.Range("K134").Value = 5418.9474
$ws.Range("L134").Value = 500506740
$ws.Range("M134").Value = -2883.9474
$ws.Range("N134").Value = -500511810

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H41").Value = 5926.6665
$ws.Range("I41").Value = 5926.6665
$ws.Range("J41").Value = 0
$ws.Range("K41").Value = 5926.6665
$ws.Range("L41").Value = 0
$ws.Range("M41").Value = -5498.6665
$ws.Range("N41").ClearContents()
$ws.Range("H50").Value = 9773.333000000001
$ws.Range("J50").Value = 9773.333000000001
$ws.Range("L50").Value = 9773.333000000001
$ws.Range("N50").Value = -11023.333
$ws.Range("H51").Value = 9388.666999999999
$ws.Range("I51").Value = 9200
$ws.Range("J51").Value = 9412.25
$ws.Range("K51").Value = 9200
$ws.Range("L51").Value = 9412.25
$ws.Range("M51").Value = -8464
$ws.Range("N51").Value = -10884.25
$ws.Range("H59").Value = 13458.556
$ws.Range("J59").Value = 13458.556
$ws.Range("L59").Value = 13458.556
$ws.Range("N59").Value = -15748.556
$ws.Range("H60").Value = 9267.4375
$ws.Range("I60").Value = 0
$ws.Range("J60").Value = 9267.4375
$ws.Range("K60").Value = 0
$ws.Range("L60").Value = 9267.4375
$ws.Range("M60").ClearContents()
$ws.Range("N60").Value = -10289.4375
$ws.Range("H61").Value = 9388.666999999999
$ws.Range("I61").Value = 9200
$ws.Range("J61").Value = 9412.25
$ws.Range("K61").Value = 9200
$ws.Range("L61").Value = 9412.25
$ws.Range("M61").Value = -8852
$ws.Range("N61").Value = -10108.25
$ws.Range("H68").Value = 17581.334
$ws.Range("J68").Value = 17581.334
$ws.Range("L68").Value = 17581.334
$ws.Range("N68").Value = -19079.334
$ws.Range("H71").Value = 17581.334
$ws.Range("J71").Value = 17581.334
$ws.Range("L71").Value = 52744.00199999999
$ws.Range("N71").Value = -60232.00199999999
$ws.Range("H74").Value = 13194.875
$ws.Range("J74").Value = 13194.875
$ws.Range("L74").Value = 13194.875
$ws.Range("N74").Value = -14942.875
$ws.Range("H77").Value = 13194.875
$ws.Range("J77").Value = 13194.875
$ws.Range("L77").Value = 39584.625
$ws.Range("N77").Value = -48320.625
$ws.Range("H132").Value = 32784.156
$ws.Range("I132").Value = 48963.668
$ws.Range("J132").Value = 1896
$ws.Range("K132").Value = 146891.004
$ws.Range("L132").Value = 5688
$ws.Range("M132").Value = -144361.004
$ws.Range("N132").Value = -10748

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H98").Value = 20000
$ws.Range("J98").Value = 20000
$ws.Range("L98").Value = 20000
$ws.Range("N98").Value = -25990
$ws.Range("H126").Value = 2038
$ws.Range("I126").Value = 0
$ws.Range("J126").Value = 2038
$ws.Range("K126").Value = 0
$ws.Range("L126").Value = 6114
$ws.Range("M126").ClearContents()
$ws.Range("N126").Value = -11054

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H64").Value = 20000
$ws.Range("J64").Value = 20000
$ws.Range("L64").Value = 20000
$ws.Range("N64").Value = -20496
$ws.Range("H67").Value = 20000
$ws.Range("J67").Value = 20000
$ws.Range("L67").Value = 20000
$ws.Range("N67").Value = -21716
$ws.Range("H111").Value = 0
$ws.Range("J111").Value = 0
$ws.Range("L111").Value = 0
$ws.Range("N111").ClearContents()
$ws.Range("H132").Value = 3680.1794
$ws.Range("I132").Value = 814.4545000000001
$ws.Range("J132").Value = 7388.7646
$ws.Range("K132").Value = 2443.3635
$ws.Range("L132").Value = 22166.2938
$ws.Range("M132").Value = 86.63649999999961
$ws.Range("N132").Value = -27226.2938
$ws.Range("H136").Value = 1482386.6
$ws.Range("I136").Value = 1743718.5
$ws.Range("J136").Value = 717057.5
$ws.Range("K136").Value = 5231155.5
$ws.Range("L136").Value = 2151172.5
$ws.Range("M136").Value = -5228605.5
$ws.Range("N136").Value = -2156272.5
